# Update the attendance sheet:
#  - Change the date strings in column A from DD/MM/YYYY to DD-MM-YYYY format
#  - Update a few attendance flag values that changed between the two revisions

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new date text (dashed format)
$dates = @{
    3  = "28-07-2022"
    4  = "01-08-2022"
    5  = "04-08-2022"
    6  = "08-08-2022"
    7  = "11-08-2022"
    8  = "15-08-2022"
    9  = "18-08-2022"
    10 = "22-08-2022"
    11 = "25-08-2022"
    12 = "29-08-2022"
    13 = "01-09-2022"
    14 = "05-09-2022"
    15 = "08-09-2022"
    16 = "12-09-2022"
    17 = "15-09-2022"
    18 = "19-09-2022"
    19 = "22-09-2022"
    20 = "26-09-2022"
    21 = "29-09-2022"
}

foreach ($row in $dates.Keys) {
    $cell = $ws.Cells.Item($row, 1)
    # Force the text number format first so Excel stores the dashed string as
    # literal text rather than re-parsing it as a date value, then clear the
    # formatting again so the cell keeps using the sheet's default style
    # (matching the original, un-styled cell).
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$row]
    $cell.ClearFormats()
}

# Attendance value updates
$ws.Cells.Item(3, 4).Value = 1   # D3: 0 -> 1
$ws.Cells.Item(3, 7).Value = 1   # G3: 0 -> 1

$ws.Cells.Item(6, 4).Value = 1   # D6: 0 -> 1
$ws.Cells.Item(6, 5).Value = 1   # E6: 0 -> 1
$ws.Cells.Item(6, 8).Value = 0   # H6: 1 -> 0
